$d = $word.ActiveDocument

# --- Step 1: capitalize the leading "v" of "vfbhgfhhjhg" -> "V" ---
$capRange = $d.Range(0, 1)
$capRange.Text = "V"

# --- Step 2: split the single run "Vfbhgfhhjhg" into two runs: "V" / "fbhgfhhjhg" ---
# Splitting the paragraph right after "V" (clean new-paragraph runs have no rPr),
# then deleting the inserted paragraph mark merges the paragraphs back together
# while keeping the runs distinct (same-formatted adjacent runs would otherwise
# collapse back into a single run if built via plain text assignment).
$splitPoint = $d.Range(1, 1)
$splitPoint.InsertParagraphAfter()
$paraMark = $d.Range(1, 2)
$paraMark.Delete()

# --- Step 3: append the new paragraphs of text after the first paragraph. ---
# The "_GoBack" bookmark currently sits collapsed at the end of the first
# paragraph's text (just before its paragraph mark). Inserting a literal
# carriage-return there (rather than via InsertParagraphAfter) creates the
# new paragraph *after* the bookmark's position, effectively carrying the
# bookmark forward into the newly created (currently last/empty) paragraph.
$endOfFirstPara = $d.Paragraphs(1).Range.End
$tail = $d.Range($endOfFirstPara, $endOfFirstPara)
$newText = "`rHioshgdfhi`rGhjgh`rJgh`rJk`rHj`rKlhj`rLoç`r"
$tail.InsertBefore($newText)

# --- Step 4: fill in the final paragraph's text ("i"), which must land
# before the bookmark that now lives in that (last) paragraph. ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastInsert = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$lastInsert.InsertBefore("i")
